$d = $word.ActiveDocument

function Replace-NextOccurrence($findText, $replaceText) {
    $range = $d.Content
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)
}

Replace-NextOccurrence "2024-02-02 Friday" "2024-02-03 Saturday"

Replace-NextOccurrence "68÷6=" "40÷9="
Replace-NextOccurrence "68÷2=" "79÷8="
Replace-NextOccurrence "81÷8=" "45÷6="
Replace-NextOccurrence "76÷5=" "47÷6="
Replace-NextOccurrence "67÷6=" "20÷9="
Replace-NextOccurrence "21÷2=" "81÷9="
Replace-NextOccurrence "61÷5=" "32÷7="
Replace-NextOccurrence "88÷4=" "89÷3="
Replace-NextOccurrence "62÷8=" "81÷9="
Replace-NextOccurrence "80÷5=" "23÷4="
Replace-NextOccurrence "74÷7=" "60÷6="
Replace-NextOccurrence "88÷4=" "93÷2="
Replace-NextOccurrence "64÷3=" "77÷9="
Replace-NextOccurrence "97÷8=" "62÷9="
Replace-NextOccurrence "84÷9=" "55÷8="
Replace-NextOccurrence "54÷3=" "70÷4="
Replace-NextOccurrence "16÷8=" "48÷5="
Replace-NextOccurrence "68÷5=" "80÷6="
Replace-NextOccurrence "92÷3=" "60÷2="
Replace-NextOccurrence "79÷3=" "59÷8="
Replace-NextOccurrence "12÷8=" "48÷8="
Replace-NextOccurrence "28÷4=" "97÷5="
Replace-NextOccurrence "79÷3=" "32÷8="
Replace-NextOccurrence "90÷4=" "30÷3="
Replace-NextOccurrence "39÷5=" "29÷4="
